# Fixed minor BOM errors
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14 ("D Schottky" / "Schottky Diode") -> rename part + description
$ws.Range("A14").Value = "Diode 10TQ035"
$ws.Range("B14").Value = "Schottky Rectifier"

# Row 13 ("LED0" / "Typical INFRARED GaAs LED") -> update description,
# remove the manufacturer part number and pricing (no longer available)
$ws.Range("B13").Value = "UVC LED"
$ws.Range("C13").ClearContents()
$ws.Range("G13").ClearContents()
$ws.Range("H13").ClearContents()
